# "source and destination zone are now lists"
# FirewallPolicies!D2 (SourceZone) and F2 (DestinationZone) become comma-separated
# lists instead of single zone names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirewallPolicies")

# Order matters for how new shared-strings get appended: DestinationZone (F2)
# first, then SourceZone (D2) second.
$ws.Range("F2").Value = "ifw_interbrand_1, ifw_sslvpn_1"
$ws.Range("D2").Value = "ifw_internal_1, ifw_sslvpn_1"

# Re-fit the columns that now hold longer text.
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 25
$ws.Columns.Item(6).ColumnWidth = 30

# Reflect where the author's cursor/viewport ended up after the edit.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D7").Select() | Out-Null
